$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Partnership_Emails")

$newRow = 14
$srcRow = 13

# Duplicate the formatting (incl. the date/time number-format style) of the
# row above, then overwrite the values for the new record.
$ws.Range("A$srcRow").Copy()
$ws.Range("A$newRow").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Cells.Item($newRow, 1).Value = 45854.496805555558
$ws.Cells.Item($newRow, 2).Value = "Moris Mwai"
$ws.Cells.Item($newRow, 3).Value = "Tech-Neo GmbH"
$ws.Cells.Item($newRow, 4).Value = "Am main City, Germany"
$ws.Cells.Item($newRow, 5).Value = "DE1567890"
$ws.Cells.Item($newRow, 6).Value = "morismwai1@gmail.com"
$ws.Cells.Item($newRow, 7).Value = "Partnership Offer"
